# "tam da thay doi thong tin" — Tam updated her own contact info in the
# "Lop NM03" roster: group 4's leader name, phone number, and an email
# address (which Excel turns into a clickable mailto hyperlink).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the old leader name with the new one
$ws.Range("D22").Value = "Đỗ Hoàng Băng Tâm"

# Add the phone number
$ws.Range("E22").Value = 969382833

# Widen column F a bit so the email address/hyperlink fits
$ws.Columns.Item(6).ColumnWidth = 27.3

# Put the email address text in the cell, then turn it into a mailto
# hyperlink (this also creates the "Hyperlink" cell style/font, just like
# Excel does automatically when you type an email address).
$ws.Range("F22").Value = "bangtam.12a3.tts@gmail.com"
[void]$ws.Hyperlinks.Add($ws.Range("F22"), "mailto:bangtam.12a3.tts@gmail.com")

# Leave the view scrolled/selected where Tam was working
[void]$ws.Range("G17").Select()
